# Add 2022-Q3 data:
# 1) Insert a new "2022-Q3" row at the top of the "总计" (summary) sheet's data,
#    shifting the existing quarter rows down by one and renumbering column A.
# 2) Insert a brand new worksheet named "2022-Q3" right after "总计" (and before
#    the existing "2022-Q2" sheet) with the per-fund holding detail rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update "总计" summary sheet: insert data for 2022-Q3 as the new first row
#    (row 2) and push the previously-existing rows down by one, renumbering
#    the index column (A) sequentially as we go.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Work from the bottom up so we never clobber a row before reading it.
# Each copy brings along the cell style (s="2" on the index column) from the
# row being pushed down, matching the existing formatting pattern.
$summary.Range("A5:D5").Copy($summary.Range("A6:D6"))
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = "2021-Q1"
$summary.Cells.Item(6,3).Value = 2
$summary.Cells.Item(6,4).Value = 0.02

$summary.Range("A4:D4").Copy($summary.Range("A5:D5"))
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q2"
$summary.Cells.Item(5,3).Value = 1
$summary.Cells.Item(5,4).Value = 0.01

$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2021-Q3"
$summary.Cells.Item(4,3).Value = 2
$summary.Cells.Item(4,4).Value = 0.03

$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q2"
$summary.Cells.Item(3,3).Value = 4
$summary.Cells.Item(3,4).Value = 0.34

# Row 2 already carries the correct style from the original sheet; just
# overwrite its values with the new 2022-Q3 figures.
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 2
$summary.Cells.Item(2,4).Value = 0.08

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q3" detail sheet right before "2022-Q2".
# ---------------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($existingQ2)
$newSheet.Name = "2022-Q3"

# Header row (bold / centered style matches the rest of the workbook).
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"
$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").HorizontalAlignment = -4108
$newSheet.Range("B1:H1").VerticalAlignment = -4160

# Data rows - text-like numeric columns must stay as text (leading zeros /
# fixed decimal formatting), so force the text number format before writing.
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).NumberFormat = "@"
$newSheet.Cells.Item(2,2).Value = "015784"
$newSheet.Cells.Item(2,3).Value = "中信建投中证1000指数增强A"
$newSheet.Cells.Item(2,4).NumberFormat = "@"
$newSheet.Cells.Item(2,4).Value = "8.10"
$newSheet.Cells.Item(2,5).NumberFormat = "@"
$newSheet.Cells.Item(2,5).Value = "92.20"
$newSheet.Cells.Item(2,6).NumberFormat = "@"
$newSheet.Cells.Item(2,6).Value = "0.73"
$newSheet.Cells.Item(2,7).NumberFormat = "@"
$newSheet.Cells.Item(2,7).Value = "0.0591"
$newSheet.Cells.Item(2,8).Value = 1

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).NumberFormat = "@"
$newSheet.Cells.Item(3,2).Value = "015785"
$newSheet.Cells.Item(3,3).Value = "中信建投中证1000指数增强C"
$newSheet.Cells.Item(3,4).NumberFormat = "@"
$newSheet.Cells.Item(3,4).Value = "3.32"
$newSheet.Cells.Item(3,5).NumberFormat = "@"
$newSheet.Cells.Item(3,5).Value = "92.20"
$newSheet.Cells.Item(3,6).NumberFormat = "@"
$newSheet.Cells.Item(3,6).Value = "0.73"
$newSheet.Cells.Item(3,7).NumberFormat = "@"
$newSheet.Cells.Item(3,7).Value = "0.0242"
$newSheet.Cells.Item(3,8).Value = 1

Write-Host "2022-Q3 sheet inserted and 总计 updated"
